$d = $word.ActiveDocument
$d.Content.Find.Execute("Recursion and Dynamic Programming", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Iterators", 2)
